# PGW_TestData.xlsx - update the Qty test-data rows.
#
# Row 3, column A: value changes from the number 5 to the text "6"
# Row 4, column A: value changes from the number 10 to the text "9"
# (both cells already carry the "text" cell style used throughout column A,
# so assigning a string keeps them as text / shared-string values)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PaymentPage")

$ws.Range("A3").Value = "6"
$ws.Range("A4").Value = "9"

# Reselect the data block, matching the selection state the workbook was
# saved with (A1:C4 selected).
$ws.Range("A1:C4").Select()
